# TAaCGH manual update
#  - "parent directory" wording for the Research folder: drop the
#    "~/" home-directory shorthand in favor of a plain "/" root path,
#    and reword the opening sentence about where Research lives.
#  - "THE CODE FOLDER" heading / folder references renamed to TAaCGH.
#  - add a missing TODO paragraph after the 9_mean_diff.perm.R blurb.

$d = $word.ActiveDocument

# 1) Opening SETTING UP paragraph: no longer "a directory ... within
#    the home directory (~/Research)" but "the parent directory
#    called Research".
$d.Content.Find.Execute(
    "files under a directory called " + [char]0x201C + "Research" + [char]0x201D + " within the home directory (~/Research). Within Research there will always be three other ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "files under the parent directory called " + [char]0x201C + "Research" + [char]0x201D + ". Within Research there will always be three other ",
    2) | Out-Null

# 2) "THE CODE FOLDER" heading becomes "THE TAaCGH FOLDER"
$d.Content.Find.Execute(
    "THE CODE FOLDER",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "THE TAaCGH FOLDER",
    2) | Out-Null

# 3) "...under Research/TAaCGH..." gains a leading slash
$d.Content.Find.Execute(
    "The programs will be under Research/TAaCGH and must be run from there.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The programs will be under /Research/TAaCGH and must be run from there.",
    2) | Out-Null

# 4) Every remaining "~/Research..." path loses its leading "~" (home
#    directory) in favor of a plain "/Research..." (parent directory).
#    wdFindContinue (Wrap = 1) plus ReplaceAll (2) walks every hit.
$d.Content.Find.Execute(
    "~/Research",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "/Research",
    2) | Out-Null

# 5) New TODO paragraph, inserted right after the 9_mean_diff.perm.R
#    paragraph ("This script is still under development...").
$rng = $d.Content
$rng.Find.Execute(
    "This script is still under development. For now, it only works within R and input need to be provided manually by inspecting the output from 9_mean_diff.perm.R",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$insertPos = $rng.Start
$rng.InsertParagraphAfter()
$newPara = $d.Range($insertPos + 1, $insertPos + 1)
$newPara.Text = "TODO: need to create a warning when adding columns to phenotype file to avoid overwriting variables."
